# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets share the same row layout for the rows being updated.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 301
    4  = 51
    5  = 361
    6  = 10969
    7  = 471
    8  = 95
    9  = 9
    11 = 143
    12 = 146
    13 = 16
    18 = 305
    19 = 1144
    21 = 882
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
